$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell (untouched H-column cell) used to restore border/format
# on percentage cells after a text-forcing apostrophe-prefixed assignment,
# since Excel infers a distinct (quote-prefixed) style for "NN%"-shaped text.
$pctFormatDonor = "H3"

$ws.Range('E2').Value = '2026-02-07 01:17:43'
$ws.Range('H2').Value = "'98%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H2').PasteSpecial(-4122)
$ws.Range('N2').Value = '-1.0 °C 0:52 TU'
$ws.Range('E3').Value = '2026-02-07 01:17:45'
$ws.Range('N3').Value = '-4.6 °C 0:59 TU'
$ws.Range('O3').Value = '-4.3 °C'
$ws.Range('E4').Value = '2026-02-07 01:17:48'
$ws.Range('H4').Value = "'52%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H4').PasteSpecial(-4122)
$ws.Range('N4').Value = '11.8 °C 0:59 TU'
$ws.Range('O4').Value = '12.5 °C'
$ws.Range('E5').Value = '2026-02-07 01:17:50'
$ws.Range('H5').Value = "'69%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H5').PasteSpecial(-4122)
$ws.Range('O5').Value = '9.6 °C'
$ws.Range('E6').Value = '2026-02-07 01:17:53'
$ws.Range('H6').Value = "'56%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H6').PasteSpecial(-4122)
$ws.Range('L6').Value = '26.3 km/h - 305º 0:32 TU'
$ws.Range('N6').Value = '12.3 °C 0:59 TU'
$ws.Range('O6').Value = '12.7 °C'
$ws.Range('E7').Value = '2026-02-07 01:17:55'
$ws.Range('H7').Value = "'70%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H7').PasteSpecial(-4122)
$ws.Range('J7').Value = '1002.0 hPa'
$ws.Range('L7').Value = '38.9 km/h - 333º 0:43 TU'
$ws.Range('N7').Value = '8.3 °C 0:55 TU'
$ws.Range('O7').Value = '8.9 °C'
$ws.Range('E8').Value = '2026-02-07 01:17:58'
$ws.Range('M8').Value = '5.7 °C 0:58 TU'
$ws.Range('O8').Value = '5.2 °C'
$ws.Range('E9').Value = '2026-02-07 01:18:00'
$ws.Range('N9').Value = '2.5 °C 0:48 TU'
$ws.Range('E10').Value = '2026-02-07 01:18:02'
$ws.Range('N10').Value = '6.4 °C 0:59 TU'
$ws.Range('O10').Value = '6.9 °C'
$ws.Range('E11').Value = '2026-02-07 01:18:05'
$ws.Range('H11').Value = "'94%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H11').PasteSpecial(-4122)
$ws.Range('I11').Value = '2.5 mm'
$ws.Range('J11').Value = '1005.2 hPa'
$ws.Range('N11').Value = '1.1 °C 0:48 TU'
$ws.Range('O11').Value = '1.7 °C'
$ws.Range('E12').Value = '2026-02-07 01:18:07'
$ws.Range('N12').Value = '10.5 °C 0:58 TU'
$ws.Range('E13').Value = '2026-02-07 01:18:10'
$ws.Range('N13').Value = '6.6 °C 0:54 TU'
$ws.Range('E14').Value = '2026-02-07 01:18:12'
$ws.Range('G14').Value = '67 cm'
$ws.Range('H14').Value = "'90%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H14').PasteSpecial(-4122)
$ws.Range('N14').Value = '-6.6 °C 0:47 TU'
$ws.Range('O14').Value = '-5.7 °C'
$ws.Range('E15').Value = '2026-02-07 01:18:14'
$ws.Range('H15').Value = "'75%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H15').PasteSpecial(-4122)
$ws.Range('J15').Value = '1001.1 hPa'
$ws.Range('M15').Value = '10.2 °C 0:35 TU'
$ws.Range('O15').Value = '8.4 °C'
$ws.Range('E16').Value = '2026-02-07 01:18:17'
$ws.Range('H16').Value = "'83%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H16').PasteSpecial(-4122)
$ws.Range('N16').Value = '3.7 °C 0:56 TU'
$ws.Range('O16').Value = '4.4 °C'
$ws.Range('E17').Value = '2026-02-07 01:18:19'
$ws.Range('H17').Value = "'95%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H17').PasteSpecial(-4122)
$ws.Range('J17').Value = '1004.3 hPa'
$ws.Range('N17').Value = '3.6 °C 0:59 TU'
$ws.Range('O17').Value = '3.9 °C'
$ws.Range('E18').Value = '2026-02-07 01:18:22'
$ws.Range('N18').Value = '-6.1 °C 0:59 TU'
$ws.Range('E19').Value = '2026-02-07 01:18:24'
$ws.Range('L19').Value = '10.1 km/h - 254º 0:56 TU'
$ws.Range('E20').Value = '2026-02-07 01:18:26'
$ws.Range('I20').Value = '0.2 mm'
$ws.Range('L20').Value = '5.8 km/h - 44º 0:57 TU'
$ws.Range('M20').Value = '-4.0 °C 0:59 TU'
$ws.Range('O20').Value = '-4.1 °C'
$ws.Range('E21').Value = '2026-02-07 01:18:29'
$ws.Range('H21').Value = "'58%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H21').PasteSpecial(-4122)
$ws.Range('J21').Value = '1000.5 hPa'
$ws.Range('L21').Value = '20.5 km/h - 325º 0:38 TU'
$ws.Range('N21').Value = '10.2 °C 0:44 TU'
$ws.Range('O21').Value = '10.6 °C'
$ws.Range('E22').Value = '2026-02-07 01:18:31'
$ws.Range('L22').Value = '5.8 km/h - 241º 0:34 TU'
$ws.Range('M22').Value = '6.3 °C 0:59 TU'
$ws.Range('O22').Value = '5.9 °C'
$ws.Range('E23').Value = '2026-02-07 01:18:34'
$ws.Range('L23').Value = '10.1 km/h - 47º 0:36 TU'
$ws.Range('O23').Value = '7.8 °C'
$ws.Range('E24').Value = '2026-02-07 01:18:36'
$ws.Range('H24').Value = "'76%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H24').PasteSpecial(-4122)
$ws.Range('N24').Value = '10.3 °C 0:39 TU'
$ws.Range('O24').Value = '10.8 °C'
$ws.Range('E25').Value = '2026-02-07 01:18:39'
$ws.Range('H25').Value = "'94%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H25').PasteSpecial(-4122)
$ws.Range('I25').Value = '1.8 mm'
$ws.Range('J25').Value = '1004.7 hPa'
$ws.Range('N25').Value = '0.9 °C 0:59 TU'
$ws.Range('O25').Value = '1.5 °C'
$ws.Range('E26').Value = '2026-02-07 01:18:41'
$ws.Range('I26').Value = '0.1 mm'
$ws.Range('O26').Value = '-1.4 °C'
$ws.Range('E27').Value = '2026-02-07 01:18:43'
$ws.Range('J27').Value = '1000.8 hPa'
$ws.Range('L27').Value = '9.4 km/h - 68º 0:59 TU'
$ws.Range('M27').Value = '8.8 °C 0:55 TU'
$ws.Range('O27').Value = '8.4 °C'
$ws.Range('E28').Value = '2026-02-07 01:18:46'
$ws.Range('H28').Value = "'83%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H28').PasteSpecial(-4122)
$ws.Range('J28').Value = '1002.7 hPa'
$ws.Range('N28').Value = '4.5 °C 0:45 TU'
$ws.Range('O28').Value = '4.7 °C'
$ws.Range('E29').Value = '2026-02-07 01:18:48'
$ws.Range('H29').Value = "'52%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H29').PasteSpecial(-4122)
$ws.Range('N29').Value = '12.4 °C 0:46 TU'
$ws.Range('O29').Value = '12.8 °C'
$ws.Range('E30').Value = '2026-02-07 01:18:51'
$ws.Range('H30').Value = "'77%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H30').PasteSpecial(-4122)
$ws.Range('I30').Value = '0.1 mm'
$ws.Range('N30').Value = '-4.4 °C 0:59 TU'
$ws.Range('O30').Value = '-3.8 °C'
$ws.Range('E31').Value = '2026-02-07 01:18:53'
$ws.Range('J31').Value = '1005.4 hPa'
$ws.Range('N31').Value = '3.8 °C 0:54 TU'
$ws.Range('E32').Value = '2026-02-07 01:18:56'
$ws.Range('H32').Value = "'68%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H32').PasteSpecial(-4122)
$ws.Range('L32').Value = '25.6 km/h - 284º 0:36 TU'
$ws.Range('M32').Value = '12.4 °C 0:48 TU'
$ws.Range('O32').Value = '11.2 °C'
$ws.Range('E33').Value = '2026-02-07 01:18:58'
$ws.Range('H33').Value = "'97%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H33').PasteSpecial(-4122)
$ws.Range('M33').Value = '7.2 °C 0:43 TU'
$ws.Range('O33').Value = '7.1 °C'
$ws.Range('E34').Value = '2026-02-07 01:19:01'
$ws.Range('H34').Value = "'71%"
$ws.Range($pctFormatDonor).Copy()
$ws.Range('H34').PasteSpecial(-4122)
$ws.Range('N34').Value = '6.9 °C 0:50 TU'
$ws.Range('O34').Value = '7.8 °C'
$ws.Range('E35').Value = '2026-02-07 01:19:03'
$ws.Range('G35').Value = '200 cm'
$ws.Range('I35').Value = '0.2 mm'
$ws.Range('N35').Value = '-4.1 °C 0:30 TU'
$ws.Range('E36').Value = '2026-02-07 01:19:05'
$ws.Range('J36').Value = '1005.3 hPa'
$ws.Range('L36').Value = '6.8 km/h - 57º 0:35 TU'
$ws.Range('N36').Value = '4.5 °C 0:34 TU'
$ws.Range('O36').Value = '4.9 °C'

$excel.CutCopyMode = 0
